# Updates the crypto price/volume table to the latest scraped snapshot.
# Price (D) and Volume(1h) (E) columns are plain text in this sheet (not
# numbers), so numeric-looking price strings are written with a leading
# apostrophe to force Excel to store them as text (preserving trailing
# zeros / multi-dot thousands separators) instead of auto-converting them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.755.38"
$ws.Range("E2").Value = "  -2.12%  "

# Row 3
$ws.Range("D3").Value = "'1.798.81"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'308.80"
$ws.Range("E5").Value = "  -1.81%  "

# Row 6
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.4626"
$ws.Range("E7").Value = "  +3.66%  "

# Row 8
$ws.Range("D8").Value = "'0.3720"
$ws.Range("E8").Value = "  -1.00%  "

# Row 9
$ws.Range("D9").Value = "'0.07269"
$ws.Range("E9").Value = "  -3.57%  "

# Row 10
$ws.Range("D10").Value = "'0.8565"
$ws.Range("E10").Value = "  -3.88%  "

# Row 11
$ws.Range("D11").Value = "'20.39"
$ws.Range("E11").Value = "  -2.97%  "

# Row 12
$ws.Range("D12").Value = "'1.740.46"
$ws.Range("E12").Value = "  -5.07%  "

# Row 13
$ws.Range("D13").Value = "'5.318"
$ws.Range("E13").Value = "  -1.70%  "

# Row 14
$ws.Range("E14").Value = "  -3.68%  "

# Row 15
$ws.Range("E15").Value = "  -1.02%  "

# Row 16
$ws.Range("D16").Value = "'90.70"
$ws.Range("E16").Value = "  -3.72%  "

# Row 17
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.17%  "

# Row 18
$ws.Range("D18").Value = "'0.000008644"
$ws.Range("E18").Value = "  -1.88%  "

# Row 19
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("E20").Value = "  -4.06%  "

# Row 21
$ws.Range("D21").Value = "'26.752.40"
$ws.Range("E21").Value = "  -2.17%  "

# Row 22
$ws.Range("D22").Value = "'5.291"
$ws.Range("E22").Value = "  +0.33%  "

# Row 23
$ws.Range("E23").Value = "  -2.33%  "

# Row 24
$ws.Range("D24").Value = "'1.967.39"
$ws.Range("E24").Value = "  -4.10%  "

# Row 25
$ws.Range("D25").Value = "'1.907"
$ws.Range("E25").Value = "  -3.75%  "

# Row 26
$ws.Range("D26").Value = "'150.34"
$ws.Range("E26").Value = "  -0.75%  "

# Row 27
$ws.Range("E27").Value = "  -2.30%  "

# Row 28
$ws.Range("D28").Value = "'2.140"
$ws.Range("E28").Value = "  -8.53%  "

# Row 29
$ws.Range("D29").Value = "'5.221"
$ws.Range("E29").Value = "  -2.78%  "

# Row 30
$ws.Range("D30").Value = "'114.12"
$ws.Range("E30").Value = "  -3.12%  "

# Row 31
$ws.Range("D31").Value = "'0.08900"
$ws.Range("E31").Value = "  +0.76%  "

# Row 32
$ws.Range("D32").Value = "'0.7567"
$ws.Range("E32").Value = "  -3.46%  "

# Row 33
$ws.Range("E33").Value = "  -3.21%  "

# Row 34
$ws.Range("D34").Value = "'4.436"
$ws.Range("E34").Value = "  -1.84%  "

# Row 35
$ws.Range("D35").Value = "'2.890"
$ws.Range("E35").Value = "  -0.17%  "

# Row 36
$ws.Range("D36").Value = "'0.9998"
$ws.Range("E36").Value = "  -0.17%  "

# Row 37
$ws.Range("D37").Value = "'1.123"
$ws.Range("E37").Value = "  +1.53%  "

# Row 38
$ws.Range("E38").Value = "  -2.38%  "

# Row 39
$ws.Range("D39").Value = "'0.05205"
$ws.Range("E39").Value = "  -2.22%  "

# Row 40 (coin re-ranked; name/link/price/volume all change)
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.374"
$ws.Range("E40").Value = "  +3.88%  "

# Row 41 (coin re-ranked; name/link/price/volume all change)
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.897"
$ws.Range("E41").Value = "  +0.74%  "

# Row 42 (coin re-ranked; name/link/price/volume all change)
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.177"
$ws.Range("E42").Value = "  -2.56%  "

# Row 43
$ws.Range("D43").Value = "'0.5216"
$ws.Range("E43").Value = "  -1.74%  "

# Row 44
$ws.Range("E44").Value = "  -4.68%  "

# Row 45
$ws.Range("D45").Value = "'8.507"
$ws.Range("E45").Value = "  -2.69%  "

# Row 46
$ws.Range("D46").Value = "'0.5004"
$ws.Range("E46").Value = "  -2.67%  "

# Row 47
$ws.Range("E47").Value = "  -4.39%  "

# Row 48
$ws.Range("D48").Value = "'104.01"
$ws.Range("E48").Value = "  -1.84%  "

# Row 49
$ws.Range("D49").Value = "'0.9997"
$ws.Range("E49").Value = "  -0.17%  "

# Row 50
$ws.Range("D50").Value = "'1.649"
$ws.Range("E50").Value = "  -3.31%  "

# Row 51
$ws.Range("D51").Value = "'0.06283"
$ws.Range("E51").Value = "  -1.38%  "
